$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 87 (shifts old 87,88,89,91 down to 89,90,91,93)
$ws.Rows("87:88").Insert()

# Populate the now-blank row 86 with data (copy values from row 85, new string for B86)
$ws.Range("A86").Value = "CW3M"
$ws.Range("B86").Value = "Demo_Baseline_2010-18_C375+"
$ws.Range("C86").Value = "2010-18"
$ws.Range("A86:C86").Style = "Normal"

$ws.Range("D86").Value = 930.3491414444444
$ws.Range("E86").Value = 1890.2624783333331
$ws.Range("F86").Value = 1.0534737777777776
$ws.Range("G86").Value = 270.41205844444437
$ws.Range("H86").Value = 9.8445367777777779
$ws.Range("I86").Value = 7.3367124444444443
$ws.Range("J86").Value = 8.2027718888888881
$ws.Range("K86").Value = 668.41088177777783
$ws.Range("L86").Value = 80.17382866666668
$ws.Range("M86").Value = 1419.6366374444444
$ws.Range("N86").Value = 932.39704044444443
$ws.Range("O86").Value = 5977.3932020000002
$ws.Range("P86").Value = 27412.728515555555
$ws.Range("Q86").Value = -0.43724066666666661
$ws.Range("R86").Value = -0.00015666666666666669

Write-Output "done"
